$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row for columns D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), P (Precio $/Kg)
$data = @{
    2 = @{ D = 44848; J = 1000; K = 1500; L = 2000; M = 1750; P = 583 }
    3 = @{ D = 44910; J = 1000; K = 1800; L = 2000; M = 1900; P = 633 }
    4 = @{ D = 44827; J = 1200; K = 2000; L = 2500; M = 2250; P = 750 }
    5 = @{ D = 44881; J = 500;  K = 1900; L = 2000; M = 1950; P = 650 }
    6 = @{ D = 44685; J = 400;  K = 1500; L = 2000; M = 1750; P = 583 }
    7 = @{ D = 44883; J = 500;  K = 1800; L = 2000; M = 1900; P = 633 }
    8 = @{ D = 44911; J = 700;  K = 1800; L = 2000; M = 1900; P = 633 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("J$row").Value = $vals.J
    $ws.Range("K$row").Value = $vals.K
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("P$row").Value = $vals.P
}
